$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.354.41'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '3.388.43'
$ws.Range('E3').Value = '  +1.80%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.66'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.83'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('E9').Value = '  +8.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.589'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.87'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.39%  '
$ws.Range('E12').Value = '  +5.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '684.06'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').Value = '3.928.88'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '69.348.81'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.120'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.71%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.374.55'
$ws.Range('E18').Value = '  +1.12%  '
$ws.Range('E19').Value = '  +2.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.901'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.40'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.08'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.42%  '
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('E26').Value = '  +1.81%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.62'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.51%  '
$ws.Range('E28').Value = '  +3.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.71'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.35%  '
$ws.Range('E31').Value = '  +2.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '558.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.65'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.10%  '
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.62%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('D37').Value = '3.714.06'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('E38').Value = '  +8.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('D41').Value = '0.0₃0708'
$ws.Range('E41').Value = '  +5.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.80%  '
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('E44').Value = '  +3.12%  '
$ws.Range('E45').Value = '  -0.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('E48').Value = '  +6.08%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.64'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('E51').Value = '  -1.71%  '
